# Refresh of the forecast data: each week's row now reflects the following
# week's start date together with updated forecast figures, and the
# Summary sheet's derived statistics are updated to match.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2 (was W1 / 2025-01-26)
$ws1.Range("B2").Value = "'2025-02-02"
$ws1.Range("D2").Value = 3
$ws1.Range("E2").Value = 4
$ws1.Range("F2").Value = 5
$ws1.Range("G2").Value = 6
$ws1.Range("H2").Value = 8

# Row 3 (was W2 / 2025-02-02)
$ws1.Range("B3").Value = "'2025-02-09"
$ws1.Range("D3").Value = 3
$ws1.Range("E3").Value = 4
$ws1.Range("F3").Value = 5
$ws1.Range("G3").Value = 6
$ws1.Range("H3").Value = 8

# Row 4 (was W3 / 2025-02-09)
$ws1.Range("B4").Value = "'2025-02-16"
$ws1.Range("D4").Value = 3
$ws1.Range("E4").Value = 4
$ws1.Range("F4").Value = 5
$ws1.Range("G4").Value = 6
$ws1.Range("H4").Value = 9

# Row 5 (was W4 / 2025-02-16)
$ws1.Range("B5").Value = "'2025-02-23"
$ws1.Range("D5").Value = 3
$ws1.Range("E5").Value = 4
$ws1.Range("F5").Value = 6
$ws1.Range("G5").Value = 7
$ws1.Range("H5").Value = 10

# Row 6 (was W5 / 2025-02-23)
$ws1.Range("B6").Value = "'2025-03-02"
$ws1.Range("D6").Value = 4
$ws1.Range("E6").Value = 5
$ws1.Range("F6").Value = 6
$ws1.Range("G6").Value = 8
$ws1.Range("H6").Value = 11

# Row 7 (was W6 / 2025-03-02)
$ws1.Range("B7").Value = "'2025-03-09"
$ws1.Range("D7").Value = 3
$ws1.Range("E7").Value = 4
$ws1.Range("F7").Value = 5
$ws1.Range("G7").Value = 7
$ws1.Range("H7").Value = 10

# Row 8 (was W7 / 2025-03-09)
$ws1.Range("B8").Value = "'2025-03-16"
$ws1.Range("D8").Value = 4
$ws1.Range("E8").Value = 5
$ws1.Range("F8").Value = 6
$ws1.Range("G8").Value = 8
$ws1.Range("H8").Value = 12

# Row 9 (was W8 / 2025-03-16)
$ws1.Range("B9").Value = "'2025-03-23"
$ws1.Range("D9").Value = 3
$ws1.Range("E9").Value = 4
$ws1.Range("F9").Value = 5
$ws1.Range("G9").Value = 7
$ws1.Range("H9").Value = 11

# Row 10 (was W9 / 2025-03-23)
$ws1.Range("B10").Value = "'2025-03-30"
$ws1.Range("D10").Value = 3
$ws1.Range("E10").Value = 4
$ws1.Range("F10").Value = 5
$ws1.Range("G10").Value = 6
$ws1.Range("H10").Value = 10

# Row 11 (was W10 / 2025-03-30)
$ws1.Range("B11").Value = "'2025-04-06"
$ws1.Range("D11").Value = 3
$ws1.Range("E11").Value = 4
$ws1.Range("F11").Value = 4
$ws1.Range("G11").Value = 7
$ws1.Range("H11").Value = 10

# Row 12 (was W11 / 2025-04-06)
$ws1.Range("B12").Value = "'2025-04-13"
$ws1.Range("D12").Value = 3
$ws1.Range("E12").Value = 4
$ws1.Range("F12").Value = 5
$ws1.Range("G12").Value = 7
$ws1.Range("H12").Value = 11

# Row 13 (was W12 / 2025-04-13)
$ws1.Range("B13").Value = "'2025-04-20"
$ws1.Range("D13").Value = 3
$ws1.Range("E13").Value = 4
$ws1.Range("F13").Value = 5
$ws1.Range("G13").Value = 7
$ws1.Range("H13").Value = 11

# Row 14 (was W13 / 2025-04-20)
$ws1.Range("B14").Value = "'2025-04-27"
$ws1.Range("D14").Value = 3
$ws1.Range("E14").Value = 4
$ws1.Range("F14").Value = 5
$ws1.Range("G14").Value = 7
$ws1.Range("H14").Value = 11

# Row 15 (was W14 / 2025-04-27)
$ws1.Range("B15").Value = "'2025-05-04"
$ws1.Range("D15").Value = 3
$ws1.Range("E15").Value = 4
$ws1.Range("F15").Value = 4
$ws1.Range("G15").Value = 6
$ws1.Range("H15").Value = 10

# Row 16 (was W15 / 2025-05-04)
$ws1.Range("B16").Value = "'2025-05-11"
$ws1.Range("D16").Value = 3
$ws1.Range("E16").Value = 4
$ws1.Range("F16").Value = 5
$ws1.Range("G16").Value = 7
$ws1.Range("H16").Value = 11

# Row 17 (was W16 / 2025-05-11)
$ws1.Range("B17").Value = "'2025-05-18"
$ws1.Range("D17").Value = 3
$ws1.Range("E17").Value = 4
$ws1.Range("F17").Value = 4
$ws1.Range("G17").Value = 6
$ws1.Range("H17").Value = 10

# --- Sheet 2: "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B2").Value = "2022-12-25 to 2025-01-26"
$ws2.Range("B8").Value = "1445 units"
$ws2.Range("B9").Value = "'52"
$ws2.Range("B10").Value = "'26"
$ws2.Range("B11").Value = "'13"
$ws2.Range("B12").Value = "'4"
$ws2.Range("B13").Value = "'2025-03-02"
$ws2.Range("B14").Value = "'3"
$ws2.Range("B15").Value = "'2025-04-20"
